$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing GDP column header to lowercase "gdp"
$ws.Range("B1").Value2 = "gdp"

# Add new macro variable columns: cpi, int_rt, s_p_index
$ws.Range("C1").Value2 = "cpi"
$ws.Range("D1").Value2 = "int_rt"
$ws.Range("E1").Value2 = "s_p_index"

# Copy the Time Range row (row 2) across the new columns
$ws.Range("C2").Value2 = $ws.Range("B2").Value2
$ws.Range("D2").Value2 = $ws.Range("B2").Value2
$ws.Range("E2").Value2 = $ws.Range("B2").Value2

# Set the Description row (row 3) for the new columns
$ws.Range("C3").Value2 = "Consumer Price Index"
$ws.Range("D3").Value2 = "Interest Rate"
$ws.Range("E3").Value2 = "S & P Index"

# Copy the data values (rows 4-19) from column B into the new columns
$ws.Range("C4:C19").Value2 = $ws.Range("B4:B19").Value2
$ws.Range("D4:D19").Value2 = $ws.Range("B4:B19").Value2
$ws.Range("E4:E19").Value2 = $ws.Range("B4:B19").Value2

# Match the new column widths to column B
$ws.Range("C1:E1").ColumnWidth = $ws.Range("B1").ColumnWidth

# Update the active selection as recorded in the saved workbook
$ws.Range("G7").Select()
